$wb = $excel.ActiveWorkbook

# "zh-cn" handback status sheet
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-22 08:20:23"
$wsZh.Range("H3").Value = "2016-03-22 08:20:50"
$wsZh.Range("E5").Value = "2016-03-22 08:20:23"
$wsZh.Range("H5").Value = "2016-03-22 08:20:50"

# "de-de" handback status sheet
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-22 08:20:27"
$wsDe.Range("H3").Value = "2016-03-22 08:20:56"
$wsDe.Range("E5").Value = "2016-03-22 08:20:27"
$wsDe.Range("H5").Value = "2016-03-22 08:20:56"
